# Update EUR->ARS rate: append the 2025-10-23T15:22:48Z quote as the new
# last row (row 95) of the rate-history sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (Fecha) and B (Hora) hold date-/time-looking text such as
# "2025-10-23" and "15:22:48". Format the target cells as text first so
# Excel stores the literal strings instead of auto-converting them into
# date/time serial numbers.
$ws.Range("A95:B95").NumberFormat = "@"

$ws.Range("A95").Value = "2025-10-23"
$ws.Range("B95").Value = "15:22:48"
$ws.Range("C95").Value = "1.00 EUR = 1,842.2436"

# Drop back to the workbook's default style so the new row doesn't carry
# an explicit text-format style override (matches the plain, unstyled
# cells used by every other row in the sheet).
$ws.Range("A95:C95").Style = "Normal"
